$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25; this shifts rows 25..123 down to 26..124,
# matching the existing format of the row above it.
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new record's data.
$ws.Cells.Item(25, 1).Value = 1
$ws.Cells.Item(25, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(25, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(25, 4).Value = 44991
$ws.Cells.Item(25, 5).Value = 15
$ws.Cells.Item(25, 6).Value = 100112038
$ws.Cells.Item(25, 7).Value = "Cebollín baby"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 450
$ws.Cells.Item(25, 11).Value = 2400
$ws.Cells.Item(25, 12).Value = 2500
$ws.Cells.Item(25, 13).Value = 2444
$ws.Cells.Item(25, 14).Value = "$/paquete 2 a 2,5 kilos"
$ws.Cells.Item(25, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(25, 16).Value = 2444
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"
